$d = $word.ActiveDocument
$nl = [char]11

function Replace-Text($old, $new) {
    $r = $d.Content
    $found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "FAILED to find: $old"
        return
    }
    $r.Text = $new
}

# --- Title / author / contact ---
Replace-Text "The Quest for Existence: Unraveling Life's Enigmas" "Exploring the Secrets of Ancient Egypt: A Journey Through History, Culture, and Legacy"
Replace-Text "Lucas Freeman" "Athena Lewis"
Replace-Text "lucas.freeman@academicsolution.edu" "educationist.athena@e-academia.edu"

# --- Intro paragraph, first block (stars/cosmos/quest) ---
Replace-Text "From the birth of stars to the germination of a single seed, the universe exudes an enigma that captivates and challenges our understanding of existence" "Journey to the sands of ancient Egypt, a land steeped in mystery, captivating our imagination for centuries"
Replace-Text " As we gaze upon the cosmos, we ponder our place in the grand tapestry of reality, the workings of the physical world, and the nature of consciousness itself" " Discover the stories of powerful pharaohs, the meticulous builders of colossal pyramids, and the guardians of profound knowledge"
Replace-Text " The quest for existence is an odyssey that encompasses myriad fields, including science, philosophy, and theology" (" Immerse yourselves in the hieroglyphics adorning temple walls, revealing insights into the lives, beliefs, and aspirations of a remarkable civilization. Witness the sheer ingenuity of irrigation systems, transforming deserts into flourishing fields, and marvel at the remarkable preservation methods, revealing mummies with lifelike features, providing glimpses into their lives thousands of years ago")

# --- Second block (scientific disciplines) ---
Replace-Text "Across scientific disciplines, we seek to unveil the mysteries of life's origins and evolution" "As we delve deeper into Egyptian culture, we encounter an array of gods and goddesses, their tales woven into creation myths, shaping the very fabric of their world"
Replace-Text " From the subatomic realm of quantum mechanics to the cosmic grandeur of astrophysics, researchers strive to unravel the fundamental laws that govern reality" " Their elaborate temples, adorned with intricate reliefs and statues, stand as testaments to their enduring legacy"
Replace-Text " The search for exoplanets and the exploration of extreme environments push the boundaries of our knowledge, hinting at the immense diversity of life forms that may exist beyond Earth" " Music, art, and dance intertwine, offering vibrant expressions of their creativity and spirituality. Embarking on this journey, we not only explore the past but also gain a newfound appreciation for the foundations upon which our modern world rests"

# --- Third block (philosophers) ---
Replace-Text "Concurrently, philosophers grapple with profound questions about the meaning of existence, free will, and the nature of reality" "Throughout history, Egypt has played a pivotal role in shaping the trajectory of civilization"
Replace-Text " They delve into the depths of human consciousness, exploring the relationship between our minds and the external world" " From a vibrant hub of trade, innovation, and diplomacy, linking the ancient world, to an arena where influential dynasties rose and fell"
Replace-Text " Through rigorous argumentation and analysis, philosophers challenge our assumptions and offer novel perspectives on the nature of existence" " Their contributions to astronomy, mathematics, and medicine continue to resonate today"

# --- Remove the rest of the paragraph (Introduction Continued.. through end) and replace with one closing sentence ---
# (Anchor the search with the unchanged "today." text right before the break sequence, since a
#  search string that *starts* with a line-break character is not reliably matched.)
$old4 = "today." + $nl + $nl + "Introduction Continued:" + $nl + $nl + "Theology, too, seeks to address the enigma of existence, examining the relationship between humanity and the divine. Through religious texts, rituals, and spiritual practices, believers find solace, guidance, and meaning in an often incomprehensible world. The study of theology prompts introspection, inspiring contemplation on the purpose of life and our place in the cosmos." + $nl + $nl + "Through the lens of social sciences, scholars investigate the impact of culture, society, and politics on our perception of existence. Sociologists analyze the intricate dynamics of human interaction, while anthropologists explore the diverse ways in which cultures shape individual and collective identities. Political scientists delve into the complex interplay of power, authority, and ideology, examining how these factors influence our understanding of existence." + $nl + $nl + "Introduction Continued:" + $nl + $nl + "Artistic expressions, such as literature, music, and visual art, provide unique insights into the enigma of existence. Writers, musicians, and artists explore the depths of human emotion and experience, inviting us to contemplate the complexities of life. Their works often provoke introspection and challenge our assumptions about reality, offering alternative perspectives on the nature of existence"
$new4 = "today. The legacies left by ancient Egypt span across continents, inspiring future generations and etching their indelible mark on the global stage"
Replace-Text $old4 $new4

# --- Summary paragraph ---
Replace-Text "The quest for existence is an interdisciplinary endeavor that encompasses science, philosophy, theology, social sciences, and the arts" "Unveiling the secrets of ancient Egypt is an exploration of history, culture, and legacy"
Replace-Text " It probes the deepest questions about life's origins, the nature of reality, and the meaning of human existence" " Journey through the sands of time, discovering the stories of pharaohs, unraveling the mysteries of hieroglyphics, and marveling at the ingenuity of ancient builders"
Replace-Text " By exploring these profound enigmas, we embark on a transformative journey of self-discovery and understanding, striving to unravel the mysteries that surround us and illuminate the essence of our being" " Delve into the captivating mythology, where gods and goddesses shaped the very fabric of life. Trace the artistic expressions and spiritual beliefs through music, art, and dance. Witness Egypt's transformative role throughout history, from a bustling hub of trade and innovation to a powerful dynasty shaping the ancient world. The legacies left by ancient Egypt continue to influence our modern world, inspiring awe and leaving an indelible mark on the global landscape"

# --- Add new empty paragraph at the end of the document body ---
$cr = [char]13
$endR = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endR.Text = $cr
